$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: fix "Dijstra's" -> "Dijkstra's" and merge the three
#     runs back into a single run (matching the author's original single
#     run with rPr lang="en-US" dirty="0"). Deleting the existing text
#     collapses the run list down to one (carrying over the first run's
#     rPr, dropping the err="1" flag on the old middle run), and then we
#     fill it back in with the corrected, merged text.
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Delete()
$title.TextFrame.TextRange.Text = "Graphs $([char]0x2013) Dijkstra$([char]0x2019)s, Prim$([char]0x2019)s, Indirect Heaps"

# --- Subtitle shape: add the course/term line, a blank line, and the
#     Readings line (with "Readings" and the CLRS citation as separate
#     runs, matching the source deck).
$sub = $s.Shapes.Item(2)
$subRange = $sub.TextFrame.TextRange
$subRange.Text = "CS4102, Spring 2021" + [char]13 + [char]13 + "Readings: CLRS 23.2, 24.2, 24.3"

# Split the "Readings" / ": CLRS ..." text into two separate runs so the
# formatting boundary matches the source (first run just "Readings",
# second run starting at the colon).
$thirdPara = $sub.TextFrame.TextRange.Paragraphs(3, 1)
$readingsRun = $thirdPara.Characters(1, 8)
$readingsRun.Text = "Readings"
